$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Demodulation" step to "Complex to symbols"
$ws.Range("B12").Value = "Complex to symbols"

# Highlight the modulation (E4) and demodulation (E12) rows using the
# same "Good" (green) style already applied to the other rows in
# column E, by copying the formatting from E2 (which already uses it).
$src = $ws.Range("E2")
$src.Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Restore the selection to match the saved view state
$ws.Range("E9").Select() | Out-Null
